# Added Samples and Files Tab to all tests
#
# This workbook ("startup" sheet) drives a test harness: each row is one
# "tab" of queries (a pairing of a WebExcel query and a dbExcel/Neo4j
# StatQuery against the same pair of input files). Row 2 already has the
# "CasesTab". This change appends two more tabs - "SamplesTab" (row 3) and
# "FilesTab" (row 4) - reusing the existing StatQuery cell (column C) and
# the existing Neo4jData/WebData file names (columns D/E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$samplesQuery = @"
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
 WHERE ss.disease_subtype IN ["Adenocarcinoma"]  
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS ``Sample ID``,
            ss.study_subject_id AS ``Case ID``,
            p.program_acronym AS ``Program Code``,
            s.study_acronym AS ``Arm``,
            ss.disease_subtype AS ``Diagnosis``,
            samp.tissue_type AS ``Tissue Type``,
            samp.composition AS ``Tissue Composition``,
            samp.sample_anatomic_site AS ``Sample Anatomic Site``,
            samp.method_of_sample_procurement AS ``Sample Procurement Method``
"@

$filesQuery = @"
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
 WHERE ss.disease_subtype IN ["Adenocarcinoma"]  
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS ``File Name``,
    head(labels(samp)) AS ``Association``,
    f.file_description AS ``Description``,
    f.file_format AS ``File Format``,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS ``Program Code``,
    s.study_acronym AS ``Arm``,
    ss.study_subject_id AS ``Case ID``,
    samp.sample_id AS ``Sample ID``
    order by f.file_name
"@

# Existing values reused on the new rows (column C is the shared StatQuery
# text already used by row 2; columns D/E are the shared Neo4j/Web file
# names already used by row 2).
$statQuery = $ws.Range("C2").Value2
$neo4jFile = $ws.Range("D2").Value2
$webFile   = $ws.Range("E2").Value2

# Tab-name labels first (so new shared-string entries land in the same
# order the original author typed them: both tab names, then each query).
$ws.Range("A3").Value2 = "SamplesTab"
$ws.Range("A4").Value2 = "FilesTab"

# Row 3 - SamplesTab
$ws.Range("B3").Value2 = $samplesQuery
$ws.Range("C3").Value2 = $statQuery
$ws.Range("D3").Value2 = $neo4jFile
$ws.Range("E3").Value2 = $webFile

# Row 4 - FilesTab
$ws.Range("B4").Value2 = $filesQuery
$ws.Range("C4").Value2 = $statQuery
$ws.Range("D4").Value2 = $neo4jFile
$ws.Range("E4").Value2 = $webFile

# Match the wrap-text styling already used by the query/StatQuery columns
# on row 2.
$ws.Range("B3:C4").WrapText = $true

# Row heights grow to fit the (longer) wrapped query text, same as row 2.
$ws.Rows.Item(3).RowHeight = 345.6
$ws.Rows.Item(4).RowHeight = 409.6

# Column widths widen slightly to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 11.94140625
$ws.Columns.Item(2).ColumnWidth = 75.2734375
$ws.Columns.Item(3).ColumnWidth = 48.71484375
$ws.Columns.Item(4).ColumnWidth = 61.2734375
$ws.Columns.Item(5).ColumnWidth = 59.94140625

# Selection ends up on the newly added FilesTab query cell.
$ws.Range("B4").Select() | Out-Null
